$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Withdraw History")

# Helper cell used as a scratch area to funnel literal text through a
# formula -> copy -> paste-special-values round trip. This forces Excel to
# store the result as a genuine text (shared-string) cell instead of
# re-inferring a number/date from a plain string assignment.
$scratch = $ws.Range("Z1")

function Set-TextValue($cellAddress, $text) {
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($cellAddress).PasteSpecial(-4163)
}

Set-TextValue "A9" "200"
Set-TextValue "B9" "0900"
Set-TextValue "C9" "12/06/2021"
Set-TextValue "D9" "Lakeland, Florida"
Set-TextValue "E9" "N/A"

$scratch.ClearContents()
